# Update "2023 Punktlighet versjon3.xlsx" - fill in column J (punctuality values)
# on the "Total" sheet for the weekday rows, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Total")

$ws.Range("J7").Value  = 0.67
$ws.Range("J8").Value  = 0.56000000000000005
$ws.Range("J9").Value  = 0.56000000000000005
$ws.Range("J10").Value = 0.55000000000000004
$ws.Range("J11").Value = 0.65
$ws.Range("J12").Value = 0.65
$ws.Range("J13").Value = 0.76
$ws.Range("J14").Value = 0.45
$ws.Range("J16").Value = 0.61
$ws.Range("J17").Value = 0.81

$ws.Activate()
$ws.Range("J15").Select()
